# Weekly update: a new "Ajo" (garlic) price record was reported for the
# Terminal Hortofrutícola Agro Chillán market. It is inserted as a new
# row 96, pushing all the existing records (previously rows 96-192) down
# by one row (now rows 97-193).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96 - shifts rows 96:192 down to 97:193
# and carries the row-above formatting (date style) onto the new row.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(96, 1).Value = 7
$ws.Cells.Item(96, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(96, 3).Value = "Ñuble"
$ws.Cells.Item(96, 4).Value = 44601
$ws.Cells.Item(96, 5).Value = 16
$ws.Cells.Item(96, 6).Value = 100112003
$ws.Cells.Item(96, 7).Value = "Ajo"
$ws.Cells.Item(96, 8).Value = "Chino"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 60
$ws.Cells.Item(96, 11).Value = 19000
$ws.Cells.Item(96, 12).Value = 20000
$ws.Cells.Item(96, 13).Value = 19500
$ws.Cells.Item(96, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(96, 15).Value = "China"
$ws.Cells.Item(96, 16).Value = 1950
$ws.Cells.Item(96, 17).Value = 10
$ws.Cells.Item(96, 18).Value = "Hortaliza"
